$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column H ("property_category") before the existing date/legislator columns.
$ws.Columns.Item(8).Insert()

# Header for the new column
$ws.Cells.Item(1, 8).Value = "property_category"

# Fill "stock" as the property_category value for every data row (2-7)
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 8).Value = "stock"
}

# Fix the malformed total text for row 5 (宏達電) from "4,600，000" to a clean "4600000",
# keeping it stored as text (matches original cell type).
$g5 = $ws.Cells.Item(5, 7)
$g5.NumberFormat = "@"
$g5.Value = "4600000"
